$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")

# The "DATE" data-type entries for CreateDate / LastUpdate columns are
# being updated to "TIMESTAMP".
$ws.Range("D13").Value = "TIMESTAMP"
$ws.Range("D15").Value = "TIMESTAMP"

# Leave the cursor where the author left it after making the edit.
$ws.Range("D19").Select()
